# "Generate Report for Handoff"
#
# The localization report previously reflected a "handed back" state; this
# run regenerates it for a fresh handoff: the Status cells move from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# associated timestamp cells are refreshed to the new generation time.
# Excel's column AutoFit (triggered by the shorter new status text) also
# narrows the Status-related columns on each sheet.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 08:56:57"

# Status columns (zh-cn / de-de) narrow to fit the new, shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 08:56:53"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 08:56:57"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
